$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (Sending cluster ECs -> Target cluster MuSCs, with new TPM-derived values)
$ws.Range("D2").Value = "MuSCs"
$ws.Range("G2").Value = 82.98768099999999
$ws.Range("H2").Value = 248.963043
$ws.Range("I2").Value = 0.4489504115427952
$ws.Range("J2").Value = 0.4489504115427952
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.7115296666666667
$ws.Range("N2").Value = 2.134589
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 59.04819699936967
$ws.Range("R2").Value = 531.433772994327
$ws.Range("S2").Value = 0.4489504115427952
$ws.Range("T2").Value = 0.4489504115427952

# Update row 3 (Sending cluster ECs -> FAPs, Target cluster stays MuSCs), with new TPM-derived values
$ws.Range("A3").Value = "FAPs"
$ws.Range("G3").Value = 63.14058933333333
$ws.Range("H3").Value = 189.421768
$ws.Range("I3").Value = 0.3415807409566563
$ws.Range("J3").Value = 0.3415807409566563
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 44.92640248148356
$ws.Range("R3").Value = 404.337622333352
$ws.Range("S3").Value = 0.3415807409566563
$ws.Range("T3").Value = 0.3415807409566563

# Update row 4 (Sending cluster FAPs -> MuSCs, Target cluster ECs -> MuSCs), with new TPM-derived values
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("G4").Value = 38.719942
$ws.Range("H4").Value = 116.159826
$ws.Range("I4").Value = 0.2094688475005485
$ws.Range("J4").Value = 0.2094688475005485
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.7115296666666667
$ws.Range("N4").Value = 2.134589
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 27.55038742461267
$ws.Range("R4").Value = 247.953486821514
$ws.Range("S4").Value = 0.2094688475005485
$ws.Range("T4").Value = 0.2094688475005485

# Remove rows 5-7 (duplicate MuSCs target-cluster combinations no longer needed)
$ws.Range("A5:T7").EntireRow.Delete() | Out-Null
